# Weekly driver report update for 2025-04-21
# Updates the "Driver Summary" sheet's Bad Drivers table (rows 3-21) with
# refreshed counts/percentages (including several adapter-driver rows that
# were re-ordered because the table is sorted by Good Roaming Calculation),
# plus Total Samples refreshes for several rows in the Good Drivers table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# --- Bad Drivers table (rows 3-20) -----------------------------------
# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.3.2 (unchanged driver)
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 135
$ws.Range("D3").Value = 86.3

# Row 4: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.1.1
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.1.1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 32
$ws.Range("D4").Value = 89.7

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3 (unchanged driver)
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 753
$ws.Range("D5").Value = 92.2

# Row 6: Intel(R) Wireless-AC 9560 160MHz - 23.40.1.1 (unchanged driver)
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = 95.5

# Row 7: now Intel(R) Dual Band Wireless-AC 8265 - 20.70.26.2
$ws.Range("A7").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.26.2"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 97.7

# Row 8: Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8 (unchanged driver)
$ws.Range("D8").Value = 97.8

# Row 10: now Intel(R) Wireless-AC 9560 160MHz - 22.80.0.9
$ws.Range("A10").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.80.0.9"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = 98.2

# Row 11: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.240.0.6
$ws.Range("A11").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.240.0.6"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 23
$ws.Range("D11").Value = 98.40000000000001

# Row 12: now Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.1.1
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.1.1"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 27

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 22.190.0.4 (unchanged driver)
$ws.Range("B13").Value = 39
$ws.Range("C13").Value = 394
$ws.Range("D13").Value = 98.59999999999999

# Row 14: now Intel(R) Dual Band Wireless-AC 8265 - 20.70.21.2
$ws.Range("A14").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.21.2"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 7

# Row 15: now Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.0.7
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.0.7"
$ws.Range("B15").Value = 155
$ws.Range("C15").Value = 1392
$ws.Range("D15").Value = 98.8

# Row 16: now Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6"
$ws.Range("B16").Value = 16
$ws.Range("C16").Value = 156
$ws.Range("D16").Value = 98.8

# Row 17: now Intel(R) Wireless-AC 9560 160MHz - 22.10.0.7
$ws.Range("A17").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.10.0.7"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 98.8

# Row 18: now Intel(R) Wi-Fi 6 AX201 160MHz - 22.130.0.5
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.130.0.5"
$ws.Range("C18").Value = 5

# Row 19: now Intel(R) Wi-Fi 6 AX201 160MHz - 22.140.0.3
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.140.0.3"
$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 148

# Row 20: now Intel(R) Wi-Fi 6 AX201 160MHz - 22.60.0.6
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.60.0.6"
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 20

# Row 21: Totals
$ws.Range("B21").Value = 260
$ws.Range("C21").Value = 3266

# --- Good Drivers table: refreshed Total Samples (column B) ----------
$ws.Range("B31").Value = 56069
$ws.Range("B32").Value = 449371
$ws.Range("B36").Value = 276086
$ws.Range("B37").Value = 625298
$ws.Range("B42").Value = 453652
$ws.Range("B50").Value = 96091
$ws.Range("B53").Value = 99549
$ws.Range("B54").Value = 77999
$ws.Range("B58").Value = 175767
$ws.Range("B59").Value = 240182
$ws.Range("B67").Value = 684728
$ws.Range("B69").Value = 210188
$ws.Range("B73").Value = 308481
$ws.Range("B77").Value = 144782
$ws.Range("B80").Value = 443223
$ws.Range("B83").Value = 109665
$ws.Range("B86").Value = 62515
